$wb = $excel.ActiveWorkbook

# --- conversionAssets sheet: add two new rows (row9, row10) ---
$wsConv = $wb.Worksheets.Item("conversionAssets")

$wsConv.Cells.Item(9, 1).Value = 8
$wsConv.Cells.Item(9, 2).Value = "DH_Heat_Pump"
$wsConv.Cells.Item(9, 3).Value = "CONVERSION"
$wsConv.Cells.Item(9, 4).Value = "HEAT_PUMP_GROUND"
$wsConv.Cells.Item(9, 5).Value = 100
$wsConv.Cells.Item(9, 6).Value = 0
$wsConv.Cells.Item(9, 7).Value = 0.65

$wsConv.Cells.Item(10, 1).Value = 9
$wsConv.Cells.Item(10, 2).Value = "DH_Peak_Boiler"
$wsConv.Cells.Item(10, 3).Value = "CONVERSION"
$wsConv.Cells.Item(10, 4).Value = "BOILER"
$wsConv.Cells.Item(10, 5).Value = 300
$wsConv.Cells.Item(10, 6).Value = 297
$wsConv.Cells.Item(10, 7).Value = 0.99

$wsConv.Range("E10").Select() | Out-Null

# --- storageAssets sheet: style + value tweaks ---
$wsStor = $wb.Worksheets.Item("storageAssets")

$wsStor.Range("L3").NumberFormat = $wsStor.Range("L5").NumberFormat
$wsStor.Range("L4").NumberFormat = $wsStor.Range("L5").NumberFormat

$wsStor.Cells.Item(11, 12).Value = 1000000000

$wsStor.Range("L12").Select() | Out-Null

Write-Output "done"
